# Apply the "Continuous" Simple Plane Test Graphs.xlsx edit:
#  - Convert the Rotation X/Y/Z "Over Time" metrics (columns K, L, M) from
#    raw radians to "radians in terms of pi" by dividing the existing
#    formula by PI(), and refresh the pass-through helper formulas in
#    columns O, P, Q, S, T, U (Translation / Scale helper columns).
#  - Update the chart titles to show the units being plotted.
#  - Nudge a couple of chart positions that drifted slightly when the
#    charts were reworked.
#  - Reset the worksheet's active selection back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the underlying formulas that feed the Rotation X/Y/Z charts
#    so that the values are expressed as radians in terms of pi
#    (i.e. divide the existing angle-sum formula by PI()).
# ---------------------------------------------------------------------

$ws.Range("K4").Formula  = "=(COS(G5) - SIN(H5) + SIN(G6) + COS(H6)) / PI()"
$ws.Range("L4").Formula  = "=(COS(F4) + SIN(H4) - SIN(F6) + COS(H6)) / PI()"
$ws.Range("M4").Formula  = "=(COS(F4) - SIN(G4) + SIN(F5) + COS(G5)) / PI()"
$ws.Range("O4").Formula  = "=I4"
$ws.Range("P4").Formula  = "=I5"
$ws.Range("Q4").Formula  = "=I6"
$ws.Range("S4").Formula  = "=F4"
$ws.Range("T4").Formula  = "=G5"
$ws.Range("U4").Formula  = "=H6"

$ws.Range("K9").Formula  = "=(COS(G10) - SIN(H10) + SIN(G11) + COS(H11)) / PI()"
$ws.Range("L9").Formula  = "=(COS(F9) + SIN(H9) - SIN(F11) + COS(H11)) / PI()"
$ws.Range("M9").Formula  = "=(COS(F9) - SIN(G9) + SIN(F10) + COS(G10)) / PI()"
$ws.Range("O9").Formula  = "=I9"
$ws.Range("P9").Formula  = "=I10"
$ws.Range("Q9").Formula  = "=I11"
$ws.Range("S9").Formula  = "=F9"
$ws.Range("T9").Formula  = "=G10"
$ws.Range("U9").Formula  = "=H11"

$ws.Range("K14").Formula = "=(COS(G15) - SIN(H15) + SIN(G16) + COS(H16)) / PI()"
$ws.Range("L14").Formula = "=(COS(F14) + SIN(H14) - SIN(F16) + COS(H16)) / PI()"
$ws.Range("M14").Formula = "=(COS(F14) - SIN(G14) + SIN(F15) + COS(G15)) / PI()"
$ws.Range("O14").Formula = "=I14"
$ws.Range("P14").Formula = "=I15"
$ws.Range("Q14").Formula = "=I16"
$ws.Range("S14").Formula = "=F14"
$ws.Range("T14").Formula = "=G15"
$ws.Range("U14").Formula = "=H16"

$ws.Range("K19").Formula = "=(COS(G20) - SIN(H20) + SIN(G21) + COS(H21)) / PI()"
$ws.Range("L19").Formula = "=(COS(F19) + SIN(H19) - SIN(F21) + COS(H21)) / PI()"
$ws.Range("M19").Formula = "=(COS(F19) - SIN(G19) + SIN(F20) + COS(G20)) / PI()"
$ws.Range("O19").Formula = "=I19"
$ws.Range("P19").Formula = "=I20"
$ws.Range("Q19").Formula = "=I21"
$ws.Range("S19").Formula = "=F19"
$ws.Range("T19").Formula = "=G20"
$ws.Range("U19").Formula = "=H21"

$ws.Range("K24").Formula = "=(COS(G25) - SIN(H25) + SIN(G26) + COS(H26)) / PI()"
$ws.Range("L24").Formula = "=(COS(F24) + SIN(H24) - SIN(F26) + COS(H26)) / PI()"
$ws.Range("M24").Formula = "=(COS(F24) - SIN(G24) + SIN(F25) + COS(G25)) / PI()"
$ws.Range("O24").Formula = "=I24"
$ws.Range("P24").Formula = "=I25"
$ws.Range("Q24").Formula = "=I26"
$ws.Range("S24").Formula = "=F24"
$ws.Range("T24").Formula = "=G25"
$ws.Range("U24").Formula = "=H26"

$ws.Range("K29").Formula = "=(COS(G30) - SIN(H30) + SIN(G31) + COS(H31)) / PI()"
$ws.Range("L29").Formula = "=(COS(F29) + SIN(H29) - SIN(F31) + COS(H31)) / PI()"
$ws.Range("M29").Formula = "=(COS(F29) - SIN(G29) + SIN(F30) + COS(G30)) / PI()"
$ws.Range("O29").Formula = "=I29"
$ws.Range("P29").Formula = "=I30"
$ws.Range("Q29").Formula = "=I31"
$ws.Range("S29").Formula = "=F29"
$ws.Range("T29").Formula = "=G30"
$ws.Range("U29").Formula = "=H31"

$ws.Range("K34").Formula = "=(COS(G35) - SIN(H35) + SIN(G36) + COS(H36)) / PI()"
$ws.Range("L34").Formula = "=(COS(F34) + SIN(H34) - SIN(F36) + COS(H36)) / PI()"
$ws.Range("M34").Formula = "=(COS(F34) - SIN(G34) + SIN(F35) + COS(G35)) / PI()"
$ws.Range("O34").Formula = "=I34"
$ws.Range("P34").Formula = "=I35"
$ws.Range("Q34").Formula = "=I36"
$ws.Range("S34").Formula = "=F34"
$ws.Range("T34").Formula = "=G35"
$ws.Range("U34").Formula = "=H36"

$ws.Range("K39").Formula = "=(COS(G40) - SIN(H40) + SIN(G41) + COS(H41)) / PI()"
$ws.Range("L39").Formula = "=(COS(F39) + SIN(H39) - SIN(F41) + COS(H41)) / PI()"
$ws.Range("M39").Formula = "=(COS(F39) - SIN(G39) + SIN(F40) + COS(G40)) / PI()"
$ws.Range("O39").Formula = "=I39"
$ws.Range("P39").Formula = "=I40"
$ws.Range("Q39").Formula = "=I41"
$ws.Range("S39").Formula = "=F39"
$ws.Range("T39").Formula = "=G40"
$ws.Range("U39").Formula = "=H41"

$ws.Range("K44").Formula = "=(COS(G45) - SIN(H45) + SIN(G46) + COS(H46)) / PI()"
$ws.Range("L44").Formula = "=(COS(F44) + SIN(H44) - SIN(F46) + COS(H46)) / PI()"
$ws.Range("M44").Formula = "=(COS(F44) - SIN(G44) + SIN(F45) + COS(G45)) / PI()"
$ws.Range("O44").Formula = "=I44"
$ws.Range("P44").Formula = "=I45"
$ws.Range("Q44").Formula = "=I46"
$ws.Range("S44").Formula = "=F44"
$ws.Range("T44").Formula = "=G45"
$ws.Range("U44").Formula = "=H46"

$ws.Range("K49").Formula = "=(COS(G50) - SIN(H50) + SIN(G51) + COS(H51)) / PI()"
$ws.Range("L49").Formula = "=(COS(F49) + SIN(H49) - SIN(F51) + COS(H51)) / PI()"
$ws.Range("M49").Formula = "=(COS(F49) - SIN(G49) + SIN(F50) + COS(G50)) / PI()"
$ws.Range("O49").Formula = "=I49"
$ws.Range("P49").Formula = "=I50"
$ws.Range("Q49").Formula = "=I51"
$ws.Range("S49").Formula = "=F49"
$ws.Range("T49").Formula = "=G50"
$ws.Range("U49").Formula = "=H51"

$ws.Range("K54").Formula = "=(COS(G55) - SIN(H55) + SIN(G56) + COS(H56)) / PI()"
$ws.Range("L54").Formula = "=(COS(F54) + SIN(H54) - SIN(F56) + COS(H56)) / PI()"
$ws.Range("M54").Formula = "=(COS(F54) - SIN(G54) + SIN(F55) + COS(G55)) / PI()"
$ws.Range("O54").Formula = "=I54"
$ws.Range("P54").Formula = "=I55"
$ws.Range("Q54").Formula = "=I56"
$ws.Range("S54").Formula = "=F54"
$ws.Range("T54").Formula = "=G55"
$ws.Range("U54").Formula = "=H56"

# ---------------------------------------------------------------------
# 2. Update the chart titles to show the units of the plotted metric.
# ---------------------------------------------------------------------

$chartObjs = $ws.ChartObjects()

$chartObjs.Item(1).Chart.ChartTitle.Text  = "Rotation X Over Time (radians in terms of pi)"
$chartObjs.Item(2).Chart.ChartTitle.Text  = "Rotation Y Over Time (radians in terms of pi)"
$chartObjs.Item(3).Chart.ChartTitle.Text  = "Rotation Z Over Time (radians in terms of pi)"
$chartObjs.Item(4).Chart.ChartTitle.Text  = "Translation X Over Time (mm)"
$chartObjs.Item(5).Chart.ChartTitle.Text  = "Translation Y Over Time (mm)"
$chartObjs.Item(6).Chart.ChartTitle.Text  = "Translation Z Over Time (mm)"
$chartObjs.Item(10).Chart.ChartTitle.Text = "Error Over Time (mm)"

# ---------------------------------------------------------------------
# 3. Nudge the Rotation Y / Rotation Z charts down slightly, matching
#    the small repositioning that happened while the charts were
#    reworked.
# ---------------------------------------------------------------------

$rotYChart = $chartObjs.Item(2)
$rotYChart.Top  = $rotYChart.Top + 13.549606299212599
$rotYChart.Left = $rotYChart.Left - 0.7937007874015762

$rotZChart = $chartObjs.Item(3)
$rotZChart.Top  = $rotZChart.Top + 27.836220472440914
$rotZChart.Left = $rotZChart.Left + 2.3244094488190967

# ---------------------------------------------------------------------
# 4. Reset the active selection back to A1.
# ---------------------------------------------------------------------

$ws.Range("A1").Select()
